$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New PSSM values for cells B2:K21 (supplemental-figures update).
# Each inner array is one row (columns B..K in order) for sheet rows 2..21.
$rowsData = @(
    @(-17.6041622724692, -0.5053796666822284, -17.6041622724692, -17.6041622724692, -17.6041622724692, -17.6041622724692, -17.6041622724692, -17.6041622724692, -17.6041622724692, -17.6041622724692),
    @(-17.6041622724692, -17.6041622724692, -17.6041622724692, -17.6041622724692, -17.6041622724692, -17.6041622724692, -17.6041622724692, 0.6882357570082518, -17.6041622724692, -17.6041622724692),
    @(-17.6041622724692, -0.5681926267572933, -0.0378629089734703, -17.6041622724692, 3.97143836928749, -17.6041622724692, -17.6041622724692, -17.6041622724692, 2.725144650813581, -17.6041622724692),
    @(-17.6041622724692, 0.1356249285492369, -17.6041622724692, -17.6041622724692, -17.6041622724692, 3.282213808144201, -17.6041622724692, -17.6041622724692, -17.6041622724692, -17.6041622724692),
    @(-17.6041622724692, -17.6041622724692, -17.6041622724692, -17.6041622724692, -17.6041622724692, -17.6041622724692, -17.6041622724692, -17.6041622724692, -17.6041622724692, -17.6041622724692),
    @(3.203171421859322, -17.6041622724692, -17.6041622724692, -17.6041622724692, -17.6041622724692, -17.6041622724692, -17.6041622724692, -17.6041622724692, -17.6041622724692, -17.6041622724692),
    @(-17.6041622724692, -17.6041622724692, -17.6041622724692, 1.337644754910637, -17.6041622724692, -17.6041622724692, -17.6041622724692, -17.6041622724692, -17.6041622724692, -17.6041622724692),
    @(3.431636260393445, -17.6041622724692, -17.6041622724692, -17.6041622724692, -17.6041622724692, -17.6041622724692, -17.6041622724692, -17.6041622724692, -17.6041622724692, -17.6041622724692),
    @(-17.6041622724692, -17.6041622724692, -17.6041622724692, -17.6041622724692, -17.6041622724692, -17.6041622724692, -17.6041622724692, -0.100608245268474, -17.6041622724692, 2.817977454130017),
    @(-17.6041622724692, -17.6041622724692, -17.6041622724692, 1.994255581758279, -17.6041622724692, 1.342963810782241, -17.6041622724692, -17.6041622724692, -17.6041622724692, 1.650628773979153),
    @(-17.6041622724692, -17.6041622724692, -17.6041622724692, -17.6041622724692, -17.6041622724692, -17.6041622724692, -17.6041622724692, -17.6041622724692, -17.6041622724692, -17.6041622724692),
    @(-17.6041622724692, -17.6041622724692, -17.6041622724692, 2.111224993951421, -17.6041622724692, -17.6041622724692, -17.6041622724692, -17.6041622724692, 0.644992855688204, 2.096807372405452),
    @(-17.6041622724692, -17.6041622724692, 0.8754324965872855, -17.6041622724692, -17.6041622724692, -17.6041622724692, -17.6041622724692, -17.6041622724692, -17.6041622724692, 1.670099843852977),
    @(-17.6041622724692, -17.6041622724692, -0.5679779256175514, -17.6041622724692, -17.6041622724692, -17.6041622724692, -17.6041622724692, -17.6041622724692, -17.6041622724692, -17.6041622724692),
    @(-17.6041622724692, -17.6041622724692, -17.6041622724692, -17.6041622724692, -17.6041622724692, -17.6041622724692, -17.6041622724692, -17.6041622724692, 2.076899374724664, -17.6041622724692),
    @(-17.6041622724692, 0.06636259027094552, -0.3875963835196986, -17.6041622724692, -17.6041622724692, -17.6041622724692, -17.6041622724692, -0.2097581952347256, 2.034756151479539, -17.6041622724692),
    @(-17.6041622724692, -17.6041622724692, -17.6041622724692, -17.6041622724692, -17.6041622724692, -17.6041622724692, 4.321921215998592, -0.3802265674925946, 1.810401164745232, -17.6041622724692),
    @(-17.6041622724692, -17.6041622724692, 3.044722775726235, -17.6041622724692, -17.6041622724692, -17.6041622724692, -17.6041622724692, 0.3175682737607975, -17.6041622724692, -17.6041622724692),
    @(-17.6041622724692, 3.269489852917035, 2.906973622415716, -17.6041622724692, 2.108876805104307, -17.6041622724692, -17.6041622724692, 3.865597277979186, -17.6041622724692, 1.231743179364276),
    @(-17.6041622724692, 2.772371200186152, -17.6041622724692, 3.196583621492744, -17.6041622724692, 2.951348989484561, -17.6041622724692, -17.6041622724692, -17.6041622724692, -17.6041622724692)
)

$startRow = 2
$startCol = 2  # column B
for ($i = 0; $i -lt $rowsData.Count; $i++) {
    $rowArr = $rowsData[$i]
    for ($j = 0; $j -lt $rowArr.Count; $j++) {
        $ws.Cells.Item($startRow + $i, $startCol + $j).Value = $rowArr[$j]
    }
}
